$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix D2 quantity from 15 to 13
$ws.Range("D2").Value = 13

# Add new row 4 with a transferred product entry
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "gفحم الصفوة 1750"
$ws.Range("C4").Value = "Wholesale / جملة"
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = "Store 2 / مخزن 2"
